# Update the repeated block of model-run output rows (105-143) on the
# active sheet with refreshed values from the latest calculation engine run.
# Every row in this block historically held identical values across all
# ten scenario columns (A:J); only the values changed, not the structure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(76.30214722287239, 110.097953417861, 151.9679578377339, 238.1132357349561, 247.5301089763919, 256.9469822178277, 266.3638554680526, 275.7807287182775, 285.1976019597133, 292.6206087008836)

for ($row = 105; $row -le 143; $row++) {
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 1]
    }
}
